$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.5977047085762024
$ws.Range("B1").Value = 1.184797525405884
$ws.Range("C1").Value = 5.283000946044922
$ws.Range("D1").Value = 1.806358456611633
$ws.Range("E1").Value = 1.159831762313843
